$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$r = $ws.Range("E4")
$r.Borders.LineStyle = 1
$r.Borders.TintAndShade = 0
